# Apply cryptos price/volume refresh from the Sun May 28 13:56:28 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "311.28") need the source
# column coerced to Text first, otherwise Excel auto-converts the assignment into
# a numeric cell instead of the original text/inline-string cell.
$textForceCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D18",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Write every updated Price (D) / Volume(1h) (E) cell.
$ws.Range("D2").Value = '27.439.57'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '1.861.25'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '311.28'
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4774'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '0.3809'
$ws.Range("E8").Value = '  +3.64%  '
$ws.Range("D9").Value = '0.07317'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D10").Value = '0.9319'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '20.77'
$ws.Range("E11").Value = '  +5.00%  '
$ws.Range("D12").Value = '0.07793'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.864.80'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '5.443'
$ws.Range("E14").Value = '  +1.87%  '
$ws.Range("D15").Value = '6.553'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = '90.15'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '0.000008812'
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '27.478.52'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").Value = '14.64'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").Value = '5.098'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("D24").Value = '1.940'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").Value = '155.30'
$ws.Range("E25").Value = '  +1.81%  '
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("D27").Value = '2.011'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '115.42'
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("D29").Value = '4.954'
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").Value = '0.08898'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = '3.330'
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '1.204'
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("D33").Value = '0.7547'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").Value = '2.708'
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '0.02041'
$ws.Range("E37").Value = '  +4.19%  '
$ws.Range("D38").Value = '0.5563'
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("D39").Value = '0.05274'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").Value = '2.986'
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").Value = '7.046'
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '8.623'
$ws.Range("E42").Value = '  +4.52%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("E44").Value = '  +3.02%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '1.011'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '1.664'
$ws.Range("E47").Value = '  +3.68%  '
$ws.Range("D48").Value = '102.99'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").Value = '67.39'
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("D50").Value = '0.06092'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '0.9133'
$ws.Range("E51").Value = '  +2.76%  '

# Restore the default (unstyled) cell style on the coerced cells so formatting
# matches the rest of the untouched data cells.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
